# Update the "想去人数" (number of people wanting to attend) column (F)
# across the four worksheets, per the generated-output refresh captured
# in the commit.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 70
$ws.Range("F4").Value  = 446
$ws.Range("F6").Value  = 332
$ws.Range("F7").Value  = 378
$ws.Range("F8").Value  = 67
$ws.Range("F9").Value  = 68
$ws.Range("F10").Value = 25
$ws.Range("F11").Value = 655
$ws.Range("F12").Value = 1505
$ws.Range("F13").Value = 5848
$ws.Range("F14").Value = 95
$ws.Range("F15").Value = 1643
$ws.Range("F16").Value = 403
$ws.Range("F17").Value = 5542
$ws.Range("F18").Value = 98
$ws.Range("F20").Value = 140
$ws.Range("F22").Value = 1577
$ws.Range("F23").Value = 826
$ws.Range("F24").Value = 31
$ws.Range("F25").Value = 92
$ws.Range("F26").Value = 1149
$ws.Range("F27").Value = 682
$ws.Range("F28").Value = 156
$ws.Range("F29").Value = 12

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 100
$ws.Range("F5").Value = 199
$ws.Range("F8").Value = 312

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9451
$ws.Range("F5").Value = 542

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 9451
$ws.Range("F5").Value  = 70
$ws.Range("F6").Value  = 446
$ws.Range("F7").Value  = 542
$ws.Range("F8").Value  = 332
$ws.Range("F9").Value  = 378
$ws.Range("F10").Value = 67
$ws.Range("F11").Value = 68
$ws.Range("F14").Value = 655
$ws.Range("F15").Value = 1505
$ws.Range("F16").Value = 5848
$ws.Range("F17").Value = 95
$ws.Range("F18").Value = 312
$ws.Range("F19").Value = 1643
$ws.Range("F22").Value = 403
$ws.Range("F25").Value = 5542
$ws.Range("F26").Value = 98
$ws.Range("F28").Value = 140
$ws.Range("F30").Value = 1577
$ws.Range("F31").Value = 826
$ws.Range("F32").Value = 31
$ws.Range("F33").Value = 92
$ws.Range("F34").Value = 1149
$ws.Range("F35").Value = 682
$ws.Range("F36").Value = 156
$ws.Range("F40").Value = 12

$wb.Save()
